$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Admin")
$ws2 = $wb.Worksheets.Item("Scouting Admin")

# --- Move the "Phone Types" block (A31:C33 on "Scouting Admin") up to
# --- the bottom of "Admin" (new rows A18:C20), preserving formatting.

# Copy formatting for the merged section-header row.
$ws2.Range("A31:C31").Copy()
$ws1.Range("A18:C18").PasteSpecial(-4122)

# Copy formatting for the two date rows individually (keeps the existing
# shared date style instead of minting a new one).
$ws2.Range("B32").Copy()
$ws1.Range("B19").PasteSpecial(-4122)
$ws2.Range("B33").Copy()
$ws1.Range("B20").PasteSpecial(-4122)

# Bring the values across.
$ws1.Range("A18").Value = $ws2.Range("A31").Value2
$ws1.Range("A19").Value = $ws2.Range("A32").Value2
$ws1.Range("B19").Value = $ws2.Range("B32").Value2
$ws1.Range("A20").Value = $ws2.Range("A33").Value2
$ws1.Range("B20").Value = $ws2.Range("B33").Value2

# Re-merge the header row on the destination, matching the source.
$ws1.Range("A18:C18").Merge()

# Mark the old location as selected before removing it, so the saved
# sheet view still points at where the block used to live.
$ws2.Activate()
$ws2.Range("A31:C33").Select()

# Remove the now-duplicated rows from "Scouting Admin".
$ws2.Range("A31:C33").EntireRow.Delete()

# Restore the view state: "Admin" tab active with its own selection.
$ws1.Activate()
$ws1.Range("A27").Select()
